# feat: add 2022-Q3 data
#
# 1. Create the new "2022-Q3" worksheet by copying the existing "2022-Q2"
#    sheet (this keeps the same header / border / font styling and keeps
#    columns A/H as numbers and B/C/D/E/F/G as text, matching the sheets
#    around it), placed immediately before "2022-Q2" (i.e. right after the
#    "总计" sheet).
# 2. Overwrite the copied sheet's data rows with the new quarter's figures.
#    D/E/F/G hold numeric-*looking* values that must stay stored as TEXT
#    (as in every other quarter sheet). Typing a numeric-looking string
#    directly into `.Value` auto-coerces it to a number, and forcing text
#    via `NumberFormat = "@"` leaves a stray style behind on the cell. So
#    the text is staged on a scratch sheet (where picking up a style
#    doesn't matter, it gets thrown away) and then brought across with
#    Copy / PasteSpecial values-only, which preserves the stored text type
#    without carrying the style along.
# 3. Insert the new 2022-Q3 summary row into "总计", shifting the existing
#    rows down by one.
#
# NOTE: worksheet object handles in this host are index-bound rather than
# identity-bound - adding/deleting/copying a sheet can silently repoint a
# variable obtained *before* that call to a different sheet afterwards. To
# stay safe, every worksheet reference below is (re-)fetched by name right
# before it's used, and never reused across a sheet-count-changing call.

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# --- 1. stage the new quarter's text-typed numbers on a scratch sheet -----
$wb.Worksheets.Add() | Out-Null
$scratch = $wb.Worksheets.Item("Sheet1")
$scratch.Range("A1:D2").NumberFormat = "@"
$scratch.Range("A1").Value = "6.39"
$scratch.Range("B1").Value = "29.75"
$scratch.Range("C1").Value = "1.40"
$scratch.Range("D1").Value = "0.0895"
$scratch.Range("A2").Value = "0.05"
$scratch.Range("B2").Value = "29.75"
$scratch.Range("C2").Value = "1.40"
$scratch.Range("D2").Value = "0.0007"

# --- 2. duplicate "2022-Q2" -> new sheet, rename to "2022-Q3" -------------
$wb.Worksheets.Item("2022-Q2").Copy($wb.Worksheets.Item("2022-Q2")) | Out-Null
$wb.Worksheets.Item("2022-Q2 (2)").Name = "2022-Q3"

# Re-fetch both sheets fresh now that the sheet collection changed.
$scratch = $wb.Worksheets.Item("Sheet1")
$q3 = $wb.Worksheets.Item("2022-Q3")

# Chinese fund names are plain text already - direct assignment is safe.
$q3.Range("C2").Value = "富国泰享回报6个月持有期混合A"
$q3.Range("C3").Value = "富国泰享回报6个月持有期混合C"

# Real numeric rank column - direct assignment is fine.
$q3.Range("H2").Value = 4
$q3.Range("H3").Value = 4

# D/E/F/G must remain TEXT cells - bring them over via paste-values so no
# extra style is introduced on $q3.
$scratch.Range("A1:D1").Copy()
$q3.Range("D2:G2").PasteSpecial(-4163)
$scratch.Range("A2:D2").Copy()
$q3.Range("D3:G3").PasteSpecial(-4163)

# --- done with the scratch sheet - remove it -------------------------------
$wb.Worksheets.Item("Sheet1").Delete() | Out-Null

# --- 3. "总计" summary sheet: insert the 2022-Q3 row -----------------------
$zj = $wb.Worksheets.Item("总计")

# Copy A4's formatting down into the new row 5 before filling values, so the
# new last row keeps the same style as the other index cells (s="2").
$zj.Range("A4").Copy($zj.Range("A5"))
$zj.Cells.Item(5, 1).Value = 3
$zj.Cells.Item(5, 2).Value = "2021-Q2"
$zj.Cells.Item(5, 3).Value = 5
$zj.Cells.Item(5, 4).Value = 1.77

# Row 4 now holds what used to be row 3's data (2022-Q1).
$zj.Cells.Item(4, 1).Value = 2
$zj.Cells.Item(4, 2).Value = "2022-Q1"
$zj.Cells.Item(4, 3).Value = 4
$zj.Cells.Item(4, 4).Value = 0.28

# Row 3 now holds what used to be row 2's data (2022-Q2).
$zj.Cells.Item(3, 1).Value = 1
$zj.Cells.Item(3, 2).Value = "2022-Q2"
$zj.Cells.Item(3, 3).Value = 2
$zj.Cells.Item(3, 4).Value = 0.1

# Row 2 gets the new 2022-Q3 summary data.
$zj.Cells.Item(2, 1).Value = 0
$zj.Cells.Item(2, 2).Value = "2022-Q3"
$zj.Cells.Item(2, 3).Value = 2
$zj.Cells.Item(2, 4).Value = 0.09
